$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1022
$ws.Range("F6").Value = 2247
$ws.Range("F10").Value = 192
$ws.Range("F11").Value = 151
$ws.Range("F12").Value = 662
$ws.Range("F13").Value = 52
$ws.Range("F15").Value = 1312
$ws.Range("F19").Value = 247

# Sheet "演出" (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 21
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 20

# Sheet "本地生活" (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6327
$ws.Range("F3").Value = 789
$ws.Range("F4").Value = 1986
$ws.Range("F5").Value = 215

# Sheet "全部类型" (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6327
$ws.Range("F3").Value = 789
$ws.Range("F4").Value = 1986
$ws.Range("F6").Value = 215
$ws.Range("F10").Value = 21
$ws.Range("F13").Value = 1022
$ws.Range("F14").Value = 12
$ws.Range("F15").Value = 20
$ws.Range("F17").Value = 2247
$ws.Range("F24").Value = 192
$ws.Range("F26").Value = 152
$ws.Range("F27").Value = 662
$ws.Range("F28").Value = 52
$ws.Range("F31").Value = 1312
$ws.Range("F43").Value = 247
